# Fruta / hortaliza, semanal
# A new weekly price record for "Haba" at Macroferia Regional de Talca needs
# to be inserted as the new row 71, pushing all the existing records
# (old rows 71-106) down by one row (new rows 72-107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 71; this shifts rows 71:106 down to 72:107
# and automatically extends the sheet dimension from R106 to R107.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly record.
$ws.Cells.Item(71, 1).Value  = 5
$ws.Cells.Item(71, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(71, 3).Value  = "Maule"
$ws.Cells.Item(71, 4).Value  = 44875
$ws.Cells.Item(71, 5).Value  = 7
$ws.Cells.Item(71, 6).Value  = 100112026
$ws.Cells.Item(71, 7).Value  = "Haba"
$ws.Cells.Item(71, 8).Value  = "Sin especificar"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 500
$ws.Cells.Item(71, 11).Value = 7000
$ws.Cells.Item(71, 12).Value = 7000
$ws.Cells.Item(71, 13).Value = 7000
$ws.Cells.Item(71, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(71, 15).Value = "Región del Maule"
$ws.Cells.Item(71, 16).Value = 280
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
